$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.272.78'
$ws.Range('E2').Value = '  +3.18%  '
$ws.Range('D3').Value = '2.319.97'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'545.13"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.28%  '
$ws.Range('D6').Value = "'131.05"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'0.581"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.39%  '
$ws.Range('D9').Value = '2.316.74'
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('D14').Value = "'23.67"
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = '60.243.25'
$ws.Range('E15').Value = '  +3.26%  '
$ws.Range('D16').Value = '2.734.39'
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '2.314.94'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('D21').Value = "'313.69"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.76%  '
$ws.Range('D22').Value = "'6.63"
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Value = "'0.996"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.35%  '
$ws.Range('D24').Value = "'63.74"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('E25').Value = '  +1.89%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  -1.98%  '
$ws.Range('E28').Value = '  +3.27%  '
$ws.Range('D29').Value = "'173.28"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.34%  '
$ws.Range('E30').Value = '  +8.13%  '
$ws.Range('E31').Value = '  +1.27%  '
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('D33').Value = "'5.92"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('E34').Value = '  +9.75%  '
$ws.Range('E35').Value = '  -0.81%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').Value = "'17.82"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').Value = '  +2.19%  '
$ws.Range('D40').Value = "'322.29"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.22%  '
$ws.Range('E41').Value = '  -1.18%  '
$ws.Range('E42').Value = '  +1.14%  '
$ws.Range('D43').Value = "'138.31"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('E44').Value = '  +0.91%  '
$ws.Range('E45').Value = '  -1.42%  '
$ws.Range('D46').Value = "'19.17"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.54%  '
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('E48').Value = '  +0.76%  '
$ws.Range('E49').Value = '  +0.56%  '
$ws.Range('E50').Value = '  +16.53%  '
